$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item(1)
# Row 8: On the Drip | Eye Drops
$ws.Range("H8").Value = 5217.6
$ws.Range("I8").Value = 104
$ws.Range("J8").Value = 10331.2
$ws.Range("K8").Value = 312
$ws.Range("L8").Value = 30993.6
$ws.Range("M8").Value = -173
$ws.Range("N8").Value = -31271.6

# Row 28: The Writing Is Not on the Wall | Enchanted Silver Ink
$ws.Range("H28").Value = 351.375
$ws.Range("I28").Value = 219.5
$ws.Range("J28").Value = 483.25
$ws.Range("K28").Value = 219.5
$ws.Range("L28").Value = 483.25
$ws.Range("M28").Value = 265.5
$ws.Range("N28").Value = -1453.25

# Row 96: Scroll Down | Grade 1 Reisui of Intelligence
$ws.Range("H96").Value = 957.1429000000001
$ws.Range("I96").Value = 957.1429000000001
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 2871.4287
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -1498.4287
$ws.Range("N96").Value = $null

# Row 116: Growing Up | Growth Formula Kappa
$ws.Range("H116").Value = 5778.5586
$ws.Range("I116").Value = 6350.8
$ws.Range("J116").Value = 4189
$ws.Range("K116").Value = 6350.8
$ws.Range("L116").Value = 4189
$ws.Range("M116").Value = -2908.8
$ws.Range("N116").Value = -11073

# Row 132: Fast-forwarding Flora | Growth Formula Lambda
$ws.Range("H132").Value = 18183158
$ws.Range("I132").Value = 1348.9783
$ws.Range("J132").Value = 111112410
$ws.Range("K132").Value = 4046.9349
$ws.Range("L132").Value = 333337230
$ws.Range("M132").Value = -1516.9349
$ws.Range("N132").Value = -333342290


# Sheet: ARM
$ws = $wb.Worksheets.Item(2)
# Row 45: Hollow Hallmarks | Mythril Ingot
$ws.Range("H45").Value = 1513.5116
$ws.Range("I45").Value = 1330.4857
$ws.Range("J45").Value = 2314.25
$ws.Range("K45").Value = 1330.4857
$ws.Range("L45").Value = 2314.25
$ws.Range("M45").Value = -953.4857
$ws.Range("N45").Value = -3068.25

# Row 119: Trial and Error | Dwarven Mythril Chainmail of Fending
$ws.Range("H119").Value = 45000
$ws.Range("J119").Value = 45000
$ws.Range("L119").Value = 45000
$ws.Range("N119").Value = -54676


# Sheet: BSM
$ws = $wb.Worksheets.Item(3)
# Row 99: Meddle in Metal | Oroshigane Ingot
$ws.Range("H99").Value = 1090.7059
$ws.Range("I99").Value = 874.2
$ws.Range("K99").Value = 874.2
$ws.Range("M99").Value = 623.8

# Row 105: Ingot to Wing It | Molybdenum Ingot
$ws.Range("H105").Value = 1533.3334
$ws.Range("I105").Value = 2000
$ws.Range("K105").Value = 2000
$ws.Range("M105").Value = -253


# Sheet: CRP
$ws = $wb.Worksheets.Item(4)
# Row 4: A Clogful of Camaraderie | Maple Clogs
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").Value = $null

# Row 43: The Long Lance of the Law | Steel Halberd
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").Value = $null

# Row 101: Everybody's Heard about the 'Berd | Doman Steel Halberd
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").Value = $null

# Row 132: Hull Lotta Damage | Ginseng Lumber
$ws.Range("H132").Value = 2745.9524
$ws.Range("I132").Value = 1759.4615
$ws.Range("J132").Value = 4349
$ws.Range("K132").Value = 5278.3845
$ws.Range("L132").Value = 13047
$ws.Range("M132").Value = -2748.3845
$ws.Range("N132").Value = -18107

# Row 134: Wood You Be Quiet | Ceiba Lumber
$ws.Range("H134").Value = 16130537
$ws.Range("I134").Value = 1255.2693
$ws.Range("K134").Value = 3765.8079
$ws.Range("M134").Value = -1230.8079


# Sheet: CUL
$ws = $wb.Worksheets.Item(5)
# Row 9: Jack of All Plates | Jack-o'-lantern
$ws.Range("H9").Value = 250000720
$ws.Range("I9").Value = 300
$ws.Range("J9").Value = 333334200
$ws.Range("K9").Value = 900
$ws.Range("L9").Value = 1000002600
$ws.Range("M9").Value = -676
$ws.Range("N9").Value = -1000003048

# Row 13: Fishy Revelations | Braised Pipira
$ws.Range("H13").Value = 1500
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 1500
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 4500
$ws.Range("M13").Value = $null
$ws.Range("N13").Value = -4836

# Row 137: Creative Chocolate | Gateau au Chocolat
$ws.Range("H137").Value = 2719972.8
$ws.Range("I137").Value = 816.36365
$ws.Range("J137").Value = 3716996.8
$ws.Range("K137").Value = 2449.09095
$ws.Range("L137").Value = 11150990.4
$ws.Range("M137").Value = 2650.90905
$ws.Range("N137").Value = -11161190.4


# Sheet: GSM
$ws = $wb.Worksheets.Item(6)
# Row 54: Tough Job Market | Horn Staff
$ws.Range("H54").Value = 6250
$ws.Range("J54").Value = 6250
$ws.Range("L54").Value = 6250
$ws.Range("N54").Value = -7030

# Row 80: Needs More Prayerbell | Hardsilver Ingot
$ws.Range("H80").Value = 2916.25
$ws.Range("J80").Value = 2851.111
$ws.Range("L80").Value = 2851.111
$ws.Range("N80").Value = -4847.111

# Row 83: With a Noise That Reaches Heaven (L) | Hardsilver Ingot
$ws.Range("H83").Value = 2916.25
$ws.Range("J83").Value = 2851.111
$ws.Range("L83").Value = 14255.555
$ws.Range("N83").Value = -24239.555

# Row 113: Copious Crystal Cannons | Manasilver Nugget
$ws.Range("H113").Value = 2950
$ws.Range("I113").Value = 2894
$ws.Range("J113").Value = 3006
$ws.Range("K113").Value = 2894
$ws.Range("L113").Value = 3006
$ws.Range("M113").Value = -724
$ws.Range("N113").Value = -7346

# Row 121: Wrap Those Wrists | Petalite Bracelet of Fending
$ws.Range("H121").Value = 20240
$ws.Range("J121").Value = 20240
$ws.Range("L121").Value = 20240
$ws.Range("N121").Value = -23734

# Row 126: Gold Rush Order | Phrygian Gold Ingot
$ws.Range("H126").Value = 1391.6562
$ws.Range("I126").Value = 1053.45
$ws.Range("J126").Value = 1955.3334
$ws.Range("K126").Value = 3160.35
$ws.Range("L126").Value = 5866.0002
$ws.Range("M126").Value = -690.3500000000004
$ws.Range("N126").Value = -10806.0002


# Sheet: LTW
$ws = $wb.Worksheets.Item(7)
# Row 119: Fit for a Friend | Swallowskin Gloves of Fending
$ws.Range("H119").Value = 41500
$ws.Range("J119").Value = 41500
$ws.Range("L119").Value = 41500
$ws.Range("N119").Value = -51176

# Row 122: Hell on Leather | Gaja Leather
$ws.Range("H122").Value = 1599.375
$ws.Range("I122").Value = 1530
$ws.Range("K122").Value = 4590
$ws.Range("M122").Value = -2140


# Sheet: WVR
$ws = $wb.Worksheets.Item(8)
# Row 81: Where the Dragonflies, the Net Catches | Crawler Silk
$ws.Range("H81").Value = 1877.9565
$ws.Range("I81").Value = 1479.4
$ws.Range("J81").Value = 1988.6666
$ws.Range("K81").Value = 2958.8
$ws.Range("L81").Value = 3977.3332
$ws.Range("M81").Value = -1897.8
$ws.Range("N81").Value = -6099.3332

# Row 84: To Kill a Dragon on Nameday (L) | Crawler Silk
$ws.Range("H84").Value = 1877.9565
$ws.Range("I84").Value = 1479.4
$ws.Range("J84").Value = 1988.6666
$ws.Range("K84").Value = 14794
$ws.Range("L84").Value = 19886.666
$ws.Range("M84").Value = -9490
$ws.Range("N84").Value = -30494.666

# Row 100: Of Great Import | Kudzu Thread
$ws.Range("H100").Value = 723.9
$ws.Range("I100").Value = 400
$ws.Range("J100").Value = 862.7143
$ws.Range("K100").Value = 800
$ws.Range("L100").Value = 1725.4286
$ws.Range("M100").Value = -259
$ws.Range("N100").Value = -2807.4286

# Row 119: A Job Well Done | Dwarven Cotton Gaskins of Fending
$ws.Range("H119").Value = 18400
$ws.Range("J119").Value = 18400
$ws.Range("L119").Value = 18400
$ws.Range("N119").Value = -28076

# Row 136: Weaving the Envelope | Sarcenet Cloth
$ws.Range("H136").Value = 3132947.8
$ws.Range("I136").Value = 3572642
$ws.Range("J136").Value = 1667299.9
$ws.Range("K136").Value = 10717926
$ws.Range("L136").Value = 5001899.699999999
$ws.Range("M136").Value = -10715376
$ws.Range("N136").Value = -5006999.699999999

